$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.593.06"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "2.653.25"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'591.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'144.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "2.653.51"
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "'27.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "3.128.92"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "63.484.28"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "2.661.36"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").Value = "'11.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "'341.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'67.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'1.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.36%  "
$ws.Range("E26").Value = "  +5.21%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").Value = "'553.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +19.94%  "
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'7.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "'1.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.63%  "
$ws.Range("D33").Value = "'1.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").Value = "0.0₃0810"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").Value = "'174.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "'4.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.76%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'0.403"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").Value = "'19.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "'1.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("D41").Value = "'170.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.16%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'40.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "'3.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'22.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "'0.0555"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.51%  "
$ws.Range("D48").Value = "'0.0961"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").Value = "'18.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("E51").Value = "  +0.47%  "
